# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# timestamps and the "Priority" column for the handoff batch of files
# (rows 7, 8, 9, 11, 13, 14 across the Overview / zh-cn / de-de sheets).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$rows = @(7, 8, 9, 11, 13, 14)

foreach ($r in $rows) {
    # Overview: "Latest HO Xliff Generate Date" column G
    $wsOverview.Range("G$r").Value = "2016-08-22 06:21:50"

    # zh-cn: "Latest Handoff Datetime" column H, "Priority" column E
    $wsZhCn.Range("H$r").Value = "2016-08-22 06:21:45"
    $wsZhCn.Range("E$r").Value = "ht"

    # de-de: "Latest Handoff Datetime" column H, "Priority" column E
    $wsDeDe.Range("H$r").Value = "2016-08-22 06:21:50"
    $wsDeDe.Range("E$r").Value = "ht"
}
